# init semaine 3 v1.0
#
# Appends, after the final paragraph ("Debut d'association ..."), three
# new paragraphs:
#   1. an empty paragraph
#   2. "Semaine 3: 26/10 - 30/10"
#   3. an empty "Paragraphedeliste" (numbered/bulleted) list-item paragraph
#
$d = $word.ActiveDocument

# Split the document's terminal paragraph mark into two paragraphs so the
# existing last paragraph ("Debut d'association ...") keeps its text/run,
# and we get a brand-new empty paragraph at the very end of the body to
# work with.
$lastPara = $d.Paragraphs.Last
$tail = $lastPara.Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()

# Collapse to the start of that freshly inserted (still empty) paragraph.
$newLast = $d.Paragraphs.Last
$insertionPoint = $newLast.Range
$insertionPoint.Collapse(1)

# Inject the two leading paragraphs (empty paragraph + "Semaine 3" heading)
# as raw WordOpenXML right before the remaining empty paragraph, which will
# become the final numbered list-item paragraph.
$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p/>
          <w:p>
            <w:r>
              <w:t>Semaine 3: 26/10 – 30/10</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
$insertionPoint.InsertXML($xml)

# The paragraph that was split off of $lastPara ("Debut d'association ...")
# already inherited its pPr (style "Paragraphedeliste" + numPr numId=1,
# ilvl=0) from the paragraph it was split from, so it is already the
# correctly-formatted trailing empty numbered list-item paragraph - no
# further action required.
